$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 204-216 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A204:A216").NumberFormat = "@"

$ws.Cells.Item(204, 1).Value = "2026-02-04"
$ws.Cells.Item(204, 2).Value = "14:18:08"
$ws.Cells.Item(204, 3).Value = "14:00"
$ws.Cells.Item(204, 4).Value = "Bathroom"
$ws.Cells.Item(204, 5).Value = "No Motion"
$ws.Cells.Item(204, 6).Value = "Inactive"

$ws.Cells.Item(205, 1).Value = "2026-02-04"
$ws.Cells.Item(205, 2).Value = "14:18:10"
$ws.Cells.Item(205, 3).Value = "14:00"
$ws.Cells.Item(205, 4).Value = "Bathroom"
$ws.Cells.Item(205, 5).Value = "No Motion"
$ws.Cells.Item(205, 6).Value = "Inactive"

$ws.Cells.Item(206, 1).Value = "2026-02-04"
$ws.Cells.Item(206, 2).Value = "14:18:12"
$ws.Cells.Item(206, 3).Value = "14:00"
$ws.Cells.Item(206, 4).Value = "Bathroom"
$ws.Cells.Item(206, 5).Value = "Motion Detected"
$ws.Cells.Item(206, 6).Value = "Active"

$ws.Cells.Item(207, 1).Value = "2026-02-04"
$ws.Cells.Item(207, 2).Value = "14:18:20"
$ws.Cells.Item(207, 3).Value = "14:00"
$ws.Cells.Item(207, 4).Value = "Bathroom"
$ws.Cells.Item(207, 5).Value = "No Motion"
$ws.Cells.Item(207, 6).Value = "Inactive"

$ws.Cells.Item(208, 1).Value = "2026-02-04"
$ws.Cells.Item(208, 2).Value = "14:18:25"
$ws.Cells.Item(208, 3).Value = "14:00"
$ws.Cells.Item(208, 4).Value = "Bathroom"
$ws.Cells.Item(208, 5).Value = "No Motion"
$ws.Cells.Item(208, 6).Value = "Inactive"

$ws.Cells.Item(209, 1).Value = "2026-02-04"
$ws.Cells.Item(209, 2).Value = "14:18:30"
$ws.Cells.Item(209, 3).Value = "14:00"
$ws.Cells.Item(209, 4).Value = "Bathroom"
$ws.Cells.Item(209, 5).Value = "No Motion"
$ws.Cells.Item(209, 6).Value = "Inactive"

$ws.Cells.Item(210, 1).Value = "2026-02-04"
$ws.Cells.Item(210, 2).Value = "14:18:35"
$ws.Cells.Item(210, 3).Value = "14:00"
$ws.Cells.Item(210, 4).Value = "Bathroom"
$ws.Cells.Item(210, 5).Value = "No Motion"
$ws.Cells.Item(210, 6).Value = "Inactive"

$ws.Cells.Item(211, 1).Value = "2026-02-04"
$ws.Cells.Item(211, 2).Value = "14:18:37"
$ws.Cells.Item(211, 3).Value = "14:00"
$ws.Cells.Item(211, 4).Value = "Bathroom"
$ws.Cells.Item(211, 5).Value = "Motion Detected"
$ws.Cells.Item(211, 6).Value = "Active"

$ws.Cells.Item(212, 1).Value = "2026-02-04"
$ws.Cells.Item(212, 2).Value = "14:18:44"
$ws.Cells.Item(212, 3).Value = "14:00"
$ws.Cells.Item(212, 4).Value = "Bathroom"
$ws.Cells.Item(212, 5).Value = "No Motion"
$ws.Cells.Item(212, 6).Value = "Inactive"

$ws.Cells.Item(213, 1).Value = "2026-02-04"
$ws.Cells.Item(213, 2).Value = "14:18:49"
$ws.Cells.Item(213, 3).Value = "14:00"
$ws.Cells.Item(213, 4).Value = "Bathroom"
$ws.Cells.Item(213, 5).Value = "No Motion"
$ws.Cells.Item(213, 6).Value = "Inactive"

$ws.Cells.Item(214, 1).Value = "2026-02-04"
$ws.Cells.Item(214, 2).Value = "14:18:55"
$ws.Cells.Item(214, 3).Value = "14:00"
$ws.Cells.Item(214, 4).Value = "Bathroom"
$ws.Cells.Item(214, 5).Value = "No Motion"
$ws.Cells.Item(214, 6).Value = "Inactive"

$ws.Cells.Item(215, 1).Value = "2026-02-04"
$ws.Cells.Item(215, 2).Value = "14:19:00"
$ws.Cells.Item(215, 3).Value = "14:00"
$ws.Cells.Item(215, 4).Value = "Bathroom"
$ws.Cells.Item(215, 5).Value = "No Motion"
$ws.Cells.Item(215, 6).Value = "Inactive"

$ws.Cells.Item(216, 1).Value = "2026-02-04"
$ws.Cells.Item(216, 2).Value = "14:19:05"
$ws.Cells.Item(216, 3).Value = "14:00"
$ws.Cells.Item(216, 4).Value = "Bathroom"
$ws.Cells.Item(216, 5).Value = "No Motion"
$ws.Cells.Item(216, 6).Value = "Inactive"

# --- Humidity sheet: rows 171-181 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A171:A181").NumberFormat = "@"
$ws.Range("E171:E181").NumberFormat = "@"

$ws.Cells.Item(171, 1).Value = "2026-02-04"
$ws.Cells.Item(171, 2).Value = "14:18:07"
$ws.Cells.Item(171, 3).Value = "14:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "76.9%"
$ws.Cells.Item(171, 6).Value = "Active"

$ws.Cells.Item(172, 1).Value = "2026-02-04"
$ws.Cells.Item(172, 2).Value = "14:18:09"
$ws.Cells.Item(172, 3).Value = "14:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "78.1%"
$ws.Cells.Item(172, 6).Value = "Active"

$ws.Cells.Item(173, 1).Value = "2026-02-04"
$ws.Cells.Item(173, 2).Value = "14:18:13"
$ws.Cells.Item(173, 3).Value = "14:00"
$ws.Cells.Item(173, 4).Value = "Bathroom"
$ws.Cells.Item(173, 5).Value = "77.3%"
$ws.Cells.Item(173, 6).Value = "Active"

$ws.Cells.Item(174, 1).Value = "2026-02-04"
$ws.Cells.Item(174, 2).Value = "14:18:18"
$ws.Cells.Item(174, 3).Value = "14:00"
$ws.Cells.Item(174, 4).Value = "Bathroom"
$ws.Cells.Item(174, 5).Value = "78.3%"
$ws.Cells.Item(174, 6).Value = "Active"

$ws.Cells.Item(175, 1).Value = "2026-02-04"
$ws.Cells.Item(175, 2).Value = "14:18:23"
$ws.Cells.Item(175, 3).Value = "14:00"
$ws.Cells.Item(175, 4).Value = "Bathroom"
$ws.Cells.Item(175, 5).Value = "77.6%"
$ws.Cells.Item(175, 6).Value = "Active"

$ws.Cells.Item(176, 1).Value = "2026-02-04"
$ws.Cells.Item(176, 2).Value = "14:18:28"
$ws.Cells.Item(176, 3).Value = "14:00"
$ws.Cells.Item(176, 4).Value = "Bathroom"
$ws.Cells.Item(176, 5).Value = "78.7%"
$ws.Cells.Item(176, 6).Value = "Active"

$ws.Cells.Item(177, 1).Value = "2026-02-04"
$ws.Cells.Item(177, 2).Value = "14:18:33"
$ws.Cells.Item(177, 3).Value = "14:00"
$ws.Cells.Item(177, 4).Value = "Bathroom"
$ws.Cells.Item(177, 5).Value = "77.9%"
$ws.Cells.Item(177, 6).Value = "Active"

$ws.Cells.Item(178, 1).Value = "2026-02-04"
$ws.Cells.Item(178, 2).Value = "14:18:38"
$ws.Cells.Item(178, 3).Value = "14:00"
$ws.Cells.Item(178, 4).Value = "Bathroom"
$ws.Cells.Item(178, 5).Value = "78.9%"
$ws.Cells.Item(178, 6).Value = "Active"

$ws.Cells.Item(179, 1).Value = "2026-02-04"
$ws.Cells.Item(179, 2).Value = "14:18:43"
$ws.Cells.Item(179, 3).Value = "14:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "78.2%"
$ws.Cells.Item(179, 6).Value = "Active"

$ws.Cells.Item(180, 1).Value = "2026-02-04"
$ws.Cells.Item(180, 2).Value = "14:18:48"
$ws.Cells.Item(180, 3).Value = "14:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "79.2%"
$ws.Cells.Item(180, 6).Value = "Active"

$ws.Cells.Item(181, 1).Value = "2026-02-04"
$ws.Cells.Item(181, 2).Value = "14:18:53"
$ws.Cells.Item(181, 3).Value = "14:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "78.2%"
$ws.Cells.Item(181, 6).Value = "Active"

# --- Temperature sheet: rows 171-181 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A171:A181").NumberFormat = "@"

$ws.Cells.Item(171, 1).Value = "2026-02-04"
$ws.Cells.Item(171, 2).Value = "14:18:08"
$ws.Cells.Item(171, 3).Value = "14:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "24.6C"
$ws.Cells.Item(171, 6).Value = "Active"

$ws.Cells.Item(172, 1).Value = "2026-02-04"
$ws.Cells.Item(172, 2).Value = "14:18:10"
$ws.Cells.Item(172, 3).Value = "14:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "24.6C"
$ws.Cells.Item(172, 6).Value = "Active"

$ws.Cells.Item(173, 1).Value = "2026-02-04"
$ws.Cells.Item(173, 2).Value = "14:18:13"
$ws.Cells.Item(173, 3).Value = "14:00"
$ws.Cells.Item(173, 4).Value = "Bathroom"
$ws.Cells.Item(173, 5).Value = "24.5C"
$ws.Cells.Item(173, 6).Value = "Active"

$ws.Cells.Item(174, 1).Value = "2026-02-04"
$ws.Cells.Item(174, 2).Value = "14:18:18"
$ws.Cells.Item(174, 3).Value = "14:00"
$ws.Cells.Item(174, 4).Value = "Bathroom"
$ws.Cells.Item(174, 5).Value = "24.5C"
$ws.Cells.Item(174, 6).Value = "Active"

$ws.Cells.Item(175, 1).Value = "2026-02-04"
$ws.Cells.Item(175, 2).Value = "14:18:23"
$ws.Cells.Item(175, 3).Value = "14:00"
$ws.Cells.Item(175, 4).Value = "Bathroom"
$ws.Cells.Item(175, 5).Value = "24.5C"
$ws.Cells.Item(175, 6).Value = "Active"

$ws.Cells.Item(176, 1).Value = "2026-02-04"
$ws.Cells.Item(176, 2).Value = "14:18:28"
$ws.Cells.Item(176, 3).Value = "14:00"
$ws.Cells.Item(176, 4).Value = "Bathroom"
$ws.Cells.Item(176, 5).Value = "24.5C"
$ws.Cells.Item(176, 6).Value = "Active"

$ws.Cells.Item(177, 1).Value = "2026-02-04"
$ws.Cells.Item(177, 2).Value = "14:18:33"
$ws.Cells.Item(177, 3).Value = "14:00"
$ws.Cells.Item(177, 4).Value = "Bathroom"
$ws.Cells.Item(177, 5).Value = "24.5C"
$ws.Cells.Item(177, 6).Value = "Active"

$ws.Cells.Item(178, 1).Value = "2026-02-04"
$ws.Cells.Item(178, 2).Value = "14:18:38"
$ws.Cells.Item(178, 3).Value = "14:00"
$ws.Cells.Item(178, 4).Value = "Bathroom"
$ws.Cells.Item(178, 5).Value = "24.5C"
$ws.Cells.Item(178, 6).Value = "Active"

$ws.Cells.Item(179, 1).Value = "2026-02-04"
$ws.Cells.Item(179, 2).Value = "14:18:44"
$ws.Cells.Item(179, 3).Value = "14:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "24.4C"
$ws.Cells.Item(179, 6).Value = "Active"

$ws.Cells.Item(180, 1).Value = "2026-02-04"
$ws.Cells.Item(180, 2).Value = "14:18:48"
$ws.Cells.Item(180, 3).Value = "14:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "24.5C"
$ws.Cells.Item(180, 6).Value = "Active"

$ws.Cells.Item(181, 1).Value = "2026-02-04"
$ws.Cells.Item(181, 2).Value = "14:18:54"
$ws.Cells.Item(181, 3).Value = "14:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "24.4C"
$ws.Cells.Item(181, 6).Value = "Active"

